$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206, pushing existing rows 206-218 down to 207-219.
$ws.Rows.Item(206).Insert()

# Populate the new row 206 with the weekly data point.
$ws.Range("A206").Value = 4
$ws.Range("B206").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C206").Value = "Los Lagos"
$ws.Range("D206").Value = 45265
$ws.Range("E206").Value = 10
$ws.Range("F206").Value = 100112052
$ws.Range("G206").Value = "Albahaca"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 90
$ws.Range("K206").Value = 7000
$ws.Range("L206").Value = 7000
$ws.Range("M206").Value = 7000
$ws.Range("N206").Value = "$/docena de matas"
$ws.Range("O206").Value = "Región Metropolitana"
$ws.Range("P206").Value = 1167
$ws.Range("Q206").Value = 6
$ws.Range("R206").Value = "Hortaliza"
